# Update "number of shares per trade" values and the values that derive from them.
$wb = $excel.ActiveWorkbook

# --- stocks sheet: shares bought and resulting value ---
$stocks = $wb.Worksheets.Item("stocks")
$stocks.Range("C2").Value = 10
$stocks.Range("D2").Value = 960.5
$stocks.Range("C3").Value = 24
$stocks.Range("D3").Value = 999.8399999999999

# --- portfolio sheet: cash and stocks totals ---
$portfolio = $wb.Worksheets.Item("portfolio")
$portfolio.Range("B2").Value = 8039.66
$portfolio.Range("B3").Value = 1960.34

# --- trades sheet: trade timestamps, shares and trade value ---
$trades = $wb.Worksheets.Item("trades")
$trades.Range("B2").Value = "01/05/2020 17:43:50"
$trades.Range("E2").Value = 10
$trades.Range("F2").Value = 960.5

$trades.Range("B3").Value = "01/05/2020 17:43:52"
$trades.Range("E3").Value = 24
$trades.Range("F3").Value = 999.8399999999999
